$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from H1 into I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I/J data values for rows 2-57
$data = @(@(9, 9), @(9, 9), @(8, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(8, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(10, 10), @(9, 9), @(8, 9), @(8, 8), @(9, 9), @(10, 10), @(10, 10), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(9, 9), @(6, 6), @(9, 9), @(7, 7), @(6, 7), @(8, 8), @(8, 8), @(7, 7), @(9, 9), @(6, 6), @(6, 6), @(9, 9), @(10, 10), @(8, 8), @(7, 7), @(6, 6), @(5, 5), @(6, 6), @(8, 8), @(7, 7), @(4, 4))

for ($k = 0; $k -lt $data.Length; $k++) {
    $row = 2 + $k
    $pair = $data[$k]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
